# Add data for 2022-04-15 (updates "through April 06" -> "through April 07")
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename sheet and update the running "through" header text/label
$ws.Name = "Through 2022-04-07"
$ws.Range("B1").Value = "April 2022 (through April 07)"

# Row 3 - Englewood
$ws.Range("B3").Value = 3

# Row 5 - Garfield Park
$ws.Range("N5").Value = 2
$ws.Range("R5").Value = 2

# Row 6 - Humboldt Park
$ws.Range("R6").Value = 2

# Row 7 - New City
$ws.Range("AD7").Value = 1

# Row 10 - Chinatown
$ws.Range("V10").Value = 1

# Row 18 - Woodlawn
$ws.Range("B18").Value = 1

# Row 29 - West Town
$ws.Range("N29").Value = 1

# Row 45 - Douglas
$ws.Range("F45").Value = 1

# Row 48 - Grand Boulevard
$ws.Range("N48").Value = 2

# Row 50 - Grand Crossing
$ws.Range("N50").Value = 1
$ws.Range("R50").Value = 1

# Row 51 - Hyde Park
$ws.Range("B51").Value = 1
$ws.Range("F51").Value = 1

# Row 57 - Albany Park
$ws.Range("B57").Value = 1

# Row 59 - Armour Square
$ws.Range("B59").Value = 1

# Row 74 - Little Italy, UIC
$ws.Range("B74").Value = 1
$ws.Range("V74").Value = 1

# Row 86 - South Chicago
$ws.Range("B86").Value = 1

# Row 92 - West Ridge
$ws.Range("J92").Value = 1
